$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Repayment Schedule")

# Insert a new (blank) column before column N to make room for a new
# "Variable Instalments" related field. This shifts the existing "Late"
# (old N) and "Over Due" (old P) columns one position to the right.
$ws.Columns("N").Insert()

# The edit leaves the "Repayment Schedule" sheet active with S5 selected.
$ws.Activate()
$ws.Range("S5").Select()
